$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full replacement data set (header row untouched).
# Columns: A Sending cluster, B Ligand symbol, C Receptor symbol, D Target cluster,
# E..T numeric metrics. 16 data rows (rows 2-17), a 4x4 cross of
# ECs/FAPs/M2/sCs sending x target clusters, Ligand=Dll1, Receptor=Notch1.
$data = @(
    @("ECs","Dll1","Notch1","ECs",3,1,13.10993833333333,39.329815,0.7359543949389752,0.7359543949389752,3,1,68.00339,204.01017,0.6265962299909886,0.6265962299909885,891.5202493576166,8023.68224421855,0.461146249314061,0.4611462493140609),
    @("ECs","Dll1","Notch1","FAPs",3,1,13.10993833333333,39.329815,0.7359543949389752,0.7359543949389752,3,1,8.99153,26.97459,0.08284967558015671,0.08284967558015671,117.8784038223167,1060.90563440085,0.06097358286248462,0.06097358286248462),
    @("ECs","Dll1","Notch1","M2",3,1,13.10993833333333,39.329815,0.7359543949389752,0.7359543949389752,3,1,10.76843933333333,32.305318,0.09922245772090688,0.09922245772090688,141.1735756062411,1270.56218045617,0.07302320383634807,0.07302320383634807),
    @("ECs","Dll1","Notch1","sCs",3,1,13.10993833333333,39.329815,0.7359543949389752,0.7359543949389752,3,1,20.764887,62.294661,0.1913316367079478,0.1913316367079478,272.226388068635,2450.037492617715,0.1408113589260815,0.1408113589260815),
    @("FAPs","Dll1","Notch1","ECs",3,1,0.5586273333333333,1.675882,0.03135973874525266,0.03135973874525266,3,1,68.00339,204.01017,0.6265962299909886,0.6265962299909885,37.98855241332666,341.8969717199399,0.01964989407127765,0.01964989407127765),
    @("FAPs","Dll1","Notch1","FAPs",3,1,0.5586273333333333,1.675882,0.03135973874525266,0.03135973874525266,3,1,8.99153,26.97459,0.08284967558015671,0.08284967558015671,5.022914426486666,45.20622983838,0.002598144181322654,0.002598144181322654),
    @("FAPs","Dll1","Notch1","M2",3,1,0.5586273333333333,1.675882,0.03135973874525266,0.03135973874525266,3,1,10.76843933333333,32.305318,0.09922245772090688,0.09922245772090688,6.015544548941777,54.139900940476,0.003111590351789518,0.003111590351789518),
    @("FAPs","Dll1","Notch1","sCs",3,1,0.5586273333333333,1.675882,0.03135973874525266,0.03135973874525266,3,1,20.764887,62.294661,0.1913316367079478,0.1913316367079478,11.599833451778,104.398501066002,0.006000110140862836,0.006000110140862836),
    @("M2","Dll1","Notch1","ECs",1,0.3333333333333333,0.09285466666666665,0.278564,0.005212595077596491,0.005212595077596491,3,1,68.00339,204.01017,0.6265962299909886,0.6265962299909885,6.314432110653332,56.82988899587999,0.003266192424091546,0.003266192424091546),
    @("M2","Dll1","Notch1","FAPs",1,0.3333333333333333,0.09285466666666665,0.278564,0.005212595077596491,0.005212595077596491,3,1,8.99153,26.97459,0.08284967558015671,0.08284967558015671,0.8349055209733331,7.514149688759999,0.0004318618111095911,0.0004318618111095911),
    @("M2","Dll1","Notch1","M2",1,0.3333333333333333,0.09285466666666665,0.278564,0.005212595077596491,0.005212595077596491,3,1,10.76843933333333,32.305318,0.09922245772090688,0.09922245772090688,0.9998998448168888,8.999098603352,0.0005172064947030252,0.0005172064947030252),
    @("M2","Dll1","Notch1","sCs",1,0.3333333333333333,0.09285466666666665,0.278564,0.005212595077596491,0.005212595077596491,3,1,20.764887,62.294661,0.1913316367079478,0.1913316367079478,1.928116660756,17.353049946804,0.0009973343476923288,0.0009973343476923288),
    @("sCs","Dll1","Notch1","ECs",3,1,4.052099666666667,12.156299,0.2274732712381757,0.2274732712381757,3,1,68.00339,204.01017,0.6265962299909886,0.6265962299909885,275.5565139512033,2480.00862556083,0.1425338941815585,0.1425338941815584),
    @("sCs","Dll1","Notch1","FAPs",3,1,4.052099666666667,12.156299,0.2274732712381757,0.2274732712381757,3,1,8.99153,26.97459,0.08284967558015671,0.08284967558015671,36.43457571582333,327.91118144241,0.01884608672523985,0.01884608672523985),
    @("sCs","Dll1","Notch1","M2",3,1,4.052099666666667,12.156299,0.2274732712381757,0.2274732712381757,3,1,10.76843933333333,32.305318,0.09922245772090688,0.09922245772090688,43.63478943312023,392.713104898082,0.02257045703806626,0.02257045703806626),
    @("sCs","Dll1","Notch1","sCs",3,1,4.052099666666667,12.156299,0.2274732712381757,0.2274732712381757,3,1,20.764887,62.294661,0.1913316367079478,0.1913316367079478,84.141391691071,757.272525219639,0.04352283329331109,0.04352283329331109)
)

$rowCount = $data.Length
$colCount = $data[0].Length

for ($r = 0; $r -lt $rowCount; $r++) {
    $rowIndex = $r + 2
    for ($c = 0; $c -lt $colCount; $c++) {
        $ws.Cells.Item($rowIndex, $c + 1).Value = $data[$r][$c]
    }
}

# Clear any stale leftover cell from the previous (shorter) table extent, if any.
$lastOldRow = 13
$lastNewRow = $rowCount + 1
if ($lastNewRow -lt $lastOldRow) {
    $clearRange = $ws.Range($ws.Cells.Item($lastNewRow + 1, 1), $ws.Cells.Item($lastOldRow, $colCount))
    $clearRange.Clear()
}
